$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates -------------------------------------------
# These cells hold plain-looking numeric text (e.g. "212.88", "0.0600").
# Assigning such text directly would make Excel auto-convert the cell to
# a real number (dropping significant trailing zeros, mangling values that
# use "." as a thousands separator, etc). To keep them as literal text we
# temporarily switch the cell to a text format, set the value, then restore
# the cell's original style so no visible formatting changes.
$dCells = @("D2", "D3", "D5", "D6", "D8", "D10", "D11", "D12", "D13", "D16", "D17", "D18", "D20", "D25", "D26", "D31", "D33", "D34", "D38", "D41", "D43", "D46", "D47", "D49")
$dStyles = @{}
foreach ($addr in $dCells) { $dStyles[$addr] = $ws.Range($addr).Style }
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "28.514.48"
$ws.Range("D3").Value = "1.586.73"
$ws.Range("D5").Value = "212.88"
$ws.Range("D6").Value = "0.493"
$ws.Range("D8").Value = "24.19"
$ws.Range("D10").Value = "0.0600"
$ws.Range("D11").Value = "0.0885"
$ws.Range("D12").Value = "1.814.38"
$ws.Range("D13").Value = "1.595.97"
$ws.Range("D16").Value = "28.528.80"
$ws.Range("D17").Value = "63.06"
$ws.Range("D18").Value = "230.89"
$ws.Range("D20").Value = "0.0$([char]0x2083)0705"
$ws.Range("D25").Value = "151.83"
$ws.Range("D26").Value = "15.21"
$ws.Range("D31").Value = "0.0469"
$ws.Range("D33").Value = "3.17"
$ws.Range("D34").Value = "1.393.19"
$ws.Range("D38").Value = "2.62"
$ws.Range("D41").Value = "0.811"
$ws.Range("D43").Value = "5.64"
$ws.Range("D46").Value = "62.89"
$ws.Range("D47").Value = "1.723.98"
$ws.Range("D49").Value = "86.92"

foreach ($addr in $dCells) { $ws.Range($addr).Style = $dStyles[$addr] }

# --- Column E (Volume/1h %) updates --------------------------------------
# These values are already padded with spaces (e.g. "  +3.73%  "), so Excel
# keeps them as plain text automatically; no format coercion to worry about.
$ws.Range("E2").Value = "  +3.73%  "
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.95%  "
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("E8").Value = "  +5.49%  "
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  -4.23%  "
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("E36").Value = "  -10.57%  "
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("E38").Value = "  +10.92%  "
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -1.02%  "
